$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.174354575127975
$ws.Range("D2").Value = 0.8631813351119915

$ws.Range("C3").Value = 1.252400066521078
$ws.Range("D3").Value = 0.2235746631123041

$ws.Range("C4").Value = 1.432129946988115
$ws.Range("D4").Value = 0.1661625180105408

$ws.Range("C5").Value = 3.522144700430246
$ws.Range("D5").Value = 0.0019193185921349

$ws.Range("C6").Value = 0.999961515074171
$ws.Range("D6").Value = 0.3282014713092709

$ws.Range("C7").Value = 1.047451143027285
$ws.Range("D7").Value = 0.3062635776570526

$ws.Range("C8").Value = 3.373737535590942
$ws.Range("D8").Value = 0.002737499145301081

$ws.Range("C9").Value = -0.0777454226573738
$ws.Range("D9").Value = 0.9387333322367319

$ws.Range("C10").Value = 2.498694673200373
$ws.Range("D10").Value = 0.02042528367103702

$ws.Range("C11").Value = 2.687317363402705
$ws.Range("D11").Value = 0.01345587153332595
